$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - update existing product
$ws.Range("A2").Value = 5151
$ws.Range("B2").Value = "ريد بل - 250 مل"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1065
$ws.Range("E2").Value = "YES"

# Row 3 - new product
$ws.Range("A3").Value = 5152
$ws.Range("B3").Value = "ريد بل فرى شوجر - 250 مل"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1065
$ws.Range("E3").Value = "YES"

# Row 4 - new product
$ws.Range("A4").Value = 5153
$ws.Range("B4").Value = "ريد بل ابيض بجوز الهند و التوت - 250 مل"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 1065
$ws.Range("E4").Value = "YES"

# Row 5 - new product
$ws.Range("A5").Value = 13928
$ws.Range("B5").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 532
$ws.Range("E5").Value = "YES"

# Row 6 - new product
$ws.Range("A6").Value = 13928
$ws.Range("B6").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 1065
$ws.Range("E6").Value = "YES"

# Row 7 - new product
$ws.Range("A7").Value = 19977
$ws.Range("B7").Value = "ريدبل كريز و توت بري  - 250 مل"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = 1050
$ws.Range("E7").Value = "YES"

# Row 8 - new product (former row 2 data moved here with updated name/tag)
$ws.Range("A8").Value = 7630
$ws.Range("B8").Value = "فيورى جولد - 400 مل"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 205
$ws.Range("E8").Value = "YES"

Write-Output "edit applied"
